$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the XT column (column B) values for rows 6-17 so the mapping
# sequence becomes 000001/000001/000001/000001, 000002 x4, 000003 x4
# instead of the previous 000002/000002/000002/000002, 000003/000004/000005/000006, 000004 x4
$ws.Range("B6").Value = "000001"
$ws.Range("B7").Value = "000001"
$ws.Range("B8").Value = "000001"
$ws.Range("B9").Value = "000001"

$ws.Range("B10").Value = "000002"
$ws.Range("B11").Value = "000002"
$ws.Range("B12").Value = "000002"
$ws.Range("B13").Value = "000002"

$ws.Range("B14").Value = "000003"
$ws.Range("B15").Value = "000003"
$ws.Range("B16").Value = "000003"
$ws.Range("B17").Value = "000003"

# Update the active selection to match the edited workbook's saved view state
$ws.Range("F10").Select()
